$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 27: C27 38 -> 39, E27 -> 0.01731793960923624
$ws.Range("C27").Value = 39
$ws.Range("E27").Value = 0.01731793960923624

# Row 34: C34 68 -> 70, E34 -> 0.03102836879432624
$ws.Range("C34").Value = 70
$ws.Range("E34").Value = 0.03102836879432624

# Row 37: C37 667 -> 671, D37 667 -> 671
$ws.Range("C37").Value = 671
$ws.Range("D37").Value = 671
